$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values per the repull / recalculation of data
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 5
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = -2
